$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint")
$lo = $ws.ListObjects.Item(1)

# Add a new row to the table for the "Manual do Usuário" task
$lo.ListRows.Add() | Out-Null
$ws.Range("A8").Value = "08/05/2024 - 15/05/2024"
$ws.Range("B8").Value = "Alterações no Manual do Usuário"
$ws.Range("C8").Value = "Bruno"
$ws.Range("D8").Value = "Pronto"

# Update existing cells to reflect renamed tasks / people
$ws.Range("C3").Value = "Cristielen"
$ws.Range("B4").Value = "Edições nos Diagramas"

$ws.Range("A9").Select()
